$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 296.35715
$ws.Range("I8").Value = 262.25
$ws.Range("K8").Value = 786.75
$ws.Range("M8").Value = -647.75

$ws.Range("H38").Value = 931.4
$ws.Range("J38").Value = 1539.6666
$ws.Range("L38").Value = 4618.9998
$ws.Range("N38").Value = -5362.9998

$ws.Range("H58").Value = 2964.6428
$ws.Range("I58").Value = 500.83334
$ws.Range("J58").Value = 4812.5
$ws.Range("K58").Value = 1502.50002
$ws.Range("L58").Value = 14437.5
$ws.Range("M58").Value = -1352.50002
$ws.Range("N58").Value = -14737.5

$ws.Range("H96").Value = 4762510
$ws.Range("I96").Value = 8928857
$ws.Range("K96").Value = 26786571
$ws.Range("M96").Value = -26785198

$ws.Range("H99").Value = 490.7
$ws.Range("J99").Value = 100
$ws.Range("L99").Value = 300
$ws.Range("N99").Value = -3296

$ws.Range("H105").Value = 88911
$ws.Range("J105").Value = 88911
$ws.Range("L105").Value = 88911
$ws.Range("N105").Value = -95899

$ws.Range("H137").Value = 27241.584
$ws.Range("J137").Value = 3256.7144
$ws.Range("L137").Value = 9770.143199999999
$ws.Range("N137").Value = -14870.1432

$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 41402.332
$ws.Range("J95").Value = 41402.332
$ws.Range("L95").Value = 41402.332
$ws.Range("N95").Value = -46894.332

$ws.Range("H102").Value = 2108.7856
$ws.Range("I102").Value = 2105.44
$ws.Range("K102").Value = 2105.44
$ws.Range("M102").Value = -483.4400000000001

$ws.Range("H104").Value = 45555
$ws.Range("J104").Value = 45555
$ws.Range("L104").Value = 45555
$ws.Range("N104").Value = -52543

$ws.Range("H132").Value = 1180.5476
$ws.Range("I132").Value = 954.4474
$ws.Range("K132").Value = 2863.3422
$ws.Range("M132").Value = -333.3422

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20408.588
$ws.Range("I20").Value = 36540.555
$ws.Range("K20").Value = 36540.555
$ws.Range("M20").Value = -36293.555

$ws.Range("H86").Value = 1570.2941
$ws.Range("I86").Value = 1346.3334
$ws.Range("K86").Value = 1346.3334
$ws.Range("M86").Value = -223.3334

$ws.Range("H89").Value = 1570.2941
$ws.Range("I89").Value = 1346.3334
$ws.Range("K89").Value = 6731.666999999999
$ws.Range("M89").Value = -1115.666999999999

$ws.Range("H99").Value = 3044.0667
$ws.Range("I99").Value = 1514.7273
$ws.Range("J99").Value = 7249.75
$ws.Range("K99").Value = 1514.7273
$ws.Range("L99").Value = 7249.75
$ws.Range("M99").Value = -16.72730000000001
$ws.Range("N99").Value = -10245.75

$ws.Range("H105").Value = 2146.4211
$ws.Range("I105").Value = 1373.1
$ws.Range("K105").Value = 1373.1
$ws.Range("M105").Value = 373.9000000000001

$ws.Range("H107").Value = 4710.8335
$ws.Range("I107").Value = 3907.45
$ws.Range("J107").Value = 6317.6
$ws.Range("K107").Value = 3907.45
$ws.Range("L107").Value = 6317.6
$ws.Range("M107").Value = -1987.45
$ws.Range("N107").Value = -10157.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4167873.5
$ws.Range("I31").Value = 4762843.5
$ws.Range("J31").Value = 3083.3333
$ws.Range("K31").Value = 4762843.5
$ws.Range("L31").Value = 3083.3333
$ws.Range("M31").Value = -4762548.5
$ws.Range("N31").Value = -3673.3333

$ws.Range("H34").Value = 4167873.5
$ws.Range("I34").Value = 4762843.5
$ws.Range("J34").Value = 3083.3333
$ws.Range("K34").Value = 4762843.5
$ws.Range("L34").Value = 3083.3333
$ws.Range("M34").Value = -4762641.5
$ws.Range("N34").Value = -3487.3333

$ws.Range("H58").Value = 14211.654
$ws.Range("I58").Value = 1455.8
$ws.Range("K58").Value = 1455.8
$ws.Range("M58").Value = -1252.8

$ws.Range("H86").Value = 34994.28
$ws.Range("I86").Value = 53825.715
$ws.Range("J86").Value = 11027
$ws.Range("K86").Value = 53825.715
$ws.Range("L86").Value = 11027
$ws.Range("M86").Value = -52702.715
$ws.Range("N86").Value = -13273

$ws.Range("H89").Value = 34994.28
$ws.Range("I89").Value = 53825.715
$ws.Range("J89").Value = 11027
$ws.Range("K89").Value = 269128.575
$ws.Range("L89").Value = 55135
$ws.Range("M89").Value = -263512.575
$ws.Range("N89").Value = -66367

$ws.Range("H105").Value = 22903.834
$ws.Range("I105").Value = 23957
$ws.Range("K105").Value = 23957
$ws.Range("M105").Value = -22210

$ws.Range("H136").Value = 14211.654
$ws.Range("I136").Value = 1455.8
$ws.Range("K136").Value = 4367.4
$ws.Range("M136").Value = -1817.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1000000000
$ws.Range("J9").Value = 1000000000
$ws.Range("L9").Value = 3000000000
$ws.Range("N9").Value = -3000000448

$ws.Range("H12").Value = 218.125
$ws.Range("J12").Value = 126.25
$ws.Range("L12").Value = 378.75
$ws.Range("N12").Value = -724.75

$ws.Range("H63").Value = 1855.5
$ws.Range("I63").Value = 1855.5
$ws.Range("K63").Value = 5566.5
$ws.Range("M63").Value = -4817.5

$ws.Range("H66").Value = 1855.5
$ws.Range("I66").Value = 1855.5
$ws.Range("K66").Value = 16699.5
$ws.Range("M66").Value = -12955.5

$ws.Range("H140").Value = 6340
$ws.Range("I140").Value = 6340
$ws.Range("K140").Value = 19020
$ws.Range("M140").Value = -13840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = ""

$ws.Range("H80").Value = 14103.363
$ws.Range("I80").Value = 7991
$ws.Range("J80").Value = 24800
$ws.Range("K80").Value = 7991
$ws.Range("L80").Value = 24800
$ws.Range("M80").Value = -6993
$ws.Range("N80").Value = -26796

$ws.Range("H83").Value = 14103.363
$ws.Range("I83").Value = 7991
$ws.Range("J83").Value = 24800
$ws.Range("K83").Value = 39955
$ws.Range("L83").Value = 124000
$ws.Range("M83").Value = -34963
$ws.Range("N83").Value = -133984

$ws.Range("H92").Value = 23350.6
$ws.Range("J92").Value = 23350.6
$ws.Range("L92").Value = 23350.6
$ws.Range("N92").Value = -27094.6

$ws.Range("H97").Value = 1694.4
$ws.Range("I97").Value = 1482
$ws.Range("K97").Value = 1482
$ws.Range("M97").Value = -986

$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3192.0625
$ws.Range("I46").Value = 2600.3333
$ws.Range("K46").Value = 2600.3333
$ws.Range("M46").Value = -2412.3333

$ws.Range("H68").Value = 3209.875
$ws.Range("I68").Value = 2963.1667
$ws.Range("K68").Value = 2963.1667
$ws.Range("M68").Value = -2214.1667

$ws.Range("H71").Value = 3209.875
$ws.Range("I71").Value = 2963.1667
$ws.Range("K71").Value = 14815.8335
$ws.Range("M71").Value = -11071.8335

$ws.Range("H100").Value = 2762.0527
$ws.Range("I100").Value = 2621.6667
$ws.Range("K100").Value = 2621.6667
$ws.Range("M100").Value = -2080.6667

$ws.Range("H106").Value = 13000
$ws.Range("J106").Value = 13000
$ws.Range("L106").Value = 13000
$ws.Range("N106").Value = -15524

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 59999.668
$ws.Range("J86").Value = 59999.668
$ws.Range("L86").Value = 59999.668
$ws.Range("N86").Value = -62245.668

$ws.Range("H89").Value = 59999.668
$ws.Range("J89").Value = 59999.668
$ws.Range("L89").Value = 299998.34
$ws.Range("N89").Value = -311230.34

$ws.Range("H132").Value = 19654.324
$ws.Range("I132").Value = 20033.611
$ws.Range("K132").Value = 60100.833
$ws.Range("M132").Value = -57570.833

$ws.Range("H136").Value = 29767.857
$ws.Range("I136").Value = 39809.133
$ws.Range("J136").Value = 4664.6665
$ws.Range("K136").Value = 119427.399
$ws.Range("L136").Value = 13993.9995
$ws.Range("M136").Value = -116877.399
$ws.Range("N136").Value = -19093.9995

Write-Host "Applied all changes"